# Auto-generated edit script: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect
# a refreshed market-data pull (scheduled runner update).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2449.3333
$ws.Range("I70").Value = 2016.8125
$ws.Range("K70").Value = 6050.4375
$ws.Range("M70").Value = -5780.4375
$ws.Range("H73").Value = 2449.3333
$ws.Range("I73").Value = 2016.8125
$ws.Range("K73").Value = 6050.4375
$ws.Range("M73").Value = -5114.4375
$ws.Range("H113").Value = 4870.788
$ws.Range("I113").Value = 3930.5
$ws.Range("J113").Value = 5755.7646
$ws.Range("K113").Value = 3930.5
$ws.Range("L113").Value = 5755.7646
$ws.Range("M113").Value = -676.5
$ws.Range("N113").Value = -12263.7646
$ws.Range("H116").Value = 3064.0527
$ws.Range("I116").Value = 2878.7693
$ws.Range("K116").Value = 2878.7693
$ws.Range("M116").Value = 563.2307000000001
$ws.Range("H125").Value = 935.8182
$ws.Range("I125").Value = 895.75
$ws.Range("J125").Value = 1042.6666
$ws.Range("K125").Value = 8061.75
$ws.Range("L125").Value = 9383.999400000001
$ws.Range("M125").Value = -5601.75
$ws.Range("N125").Value = -14303.9994
$ws.Range("H141").Value = 4680.5713
$ws.Range("I141").Value = 4680.5713
$ws.Range("K141").Value = 14041.7139
$ws.Range("M141").Value = -8861.713899999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1760.6333
$ws.Range("I61").Value = 905.7143
$ws.Range("K61").Value = 905.7143
$ws.Range("M61").Value = -693.7143
$ws.Range("H102").Value = 950
$ws.Range("I102").Value = 950
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 950
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("M102").Value = 672
$ws.Range("H124").Value = 29999.5
$ws.Range("J124").Value = 29999.5
$ws.Range("L124").Value = 29999.5
$ws.Range("N124").Value = -39819.5
$ws.Range("H125").Value = 68905
$ws.Range("J125").Value = 68905
$ws.Range("L125").Value = 68905
$ws.Range("N125").Value = -78745
$ws.Range("H132").Value = 2379.3594
$ws.Range("I132").Value = 1926.4464
$ws.Range("K132").Value = 5779.3392
$ws.Range("M132").Value = -3249.3392
$ws.Range("H136").Value = 1760.6333
$ws.Range("I136").Value = 905.7143
$ws.Range("K136").Value = 2717.1429
$ws.Range("M136").Value = -167.1428999999998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3141.3215
$ws.Range("I86").Value = 1726.5
$ws.Range("J86").Value = 4556.143
$ws.Range("K86").Value = 1726.5
$ws.Range("L86").Value = 4556.143
$ws.Range("M86").Value = -603.5
$ws.Range("N86").Value = -6802.143
$ws.Range("H89").Value = 3141.3215
$ws.Range("I89").Value = 1726.5
$ws.Range("J89").Value = 4556.143
$ws.Range("K89").Value = 8632.5
$ws.Range("L89").Value = 22780.715
$ws.Range("M89").Value = -3016.5
$ws.Range("N89").Value = -34012.715
$ws.Range("H94").Value = 3603.5334
$ws.Range("I94").Value = 2904.3
$ws.Range("K94").Value = 2904.3
$ws.Range("M94").Value = -2453.3
$ws.Range("H99").Value = 18477.041
$ws.Range("I99").Value = 22418.895
$ws.Range("J99").Value = 3498
$ws.Range("K99").Value = 22418.895
$ws.Range("L99").Value = 3498
$ws.Range("M99").Value = -20920.895
$ws.Range("N99").Value = -6494
$ws.Range("H105").Value = 3224.875
$ws.Range("I105").Value = 3224.875
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3224.875
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("M105").Value = -1477.875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1159.8572
$ws.Range("I16").Value = 853.1667
$ws.Range("K16").Value = 853.1667
$ws.Range("M16").Value = -566.1667
$ws.Range("H31").Value = 2828.9143
$ws.Range("I31").Value = 1177.5
$ws.Range("J31").Value = 5030.8
$ws.Range("K31").Value = 1177.5
$ws.Range("L31").Value = 5030.8
$ws.Range("M31").Value = -882.5
$ws.Range("N31").Value = -5620.8
$ws.Range("H34").Value = 2828.9143
$ws.Range("I34").Value = 1177.5
$ws.Range("J34").Value = 5030.8
$ws.Range("K34").Value = 1177.5
$ws.Range("L34").Value = 5030.8
$ws.Range("M34").Value = -975.5
$ws.Range("N34").Value = -5434.8
$ws.Range("H113").Value = 1159.8572
$ws.Range("I113").Value = 853.1667
$ws.Range("K113").Value = 853.1667
$ws.Range("M113").Value = 1316.8333
$ws.Range("H132").Value = 3438.9033
$ws.Range("I132").Value = 2462.56
$ws.Range("K132").Value = 7387.68
$ws.Range("M132").Value = -4857.68

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 103063.78
$ws.Range("I11").Value = 132501.72
$ws.Range("J11").Value = 31
$ws.Range("K11").Value = 397505.16
$ws.Range("L11").Value = 93
$ws.Range("M11").Value = -397365.16
$ws.Range("N11").Value = -373
$ws.Range("H37").Value = 349885.6
$ws.Range("J37").Value = 349885.6
$ws.Range("L37").Value = 1049656.8
$ws.Range("N37").Value = -1049880.8
$ws.Range("H64").Value = 2478.9
$ws.Range("I64").Value = 2389
$ws.Range("K64").Value = 7167
$ws.Range("M64").Value = -6897
$ws.Range("H67").Value = 2478.9
$ws.Range("I67").Value = 2389
$ws.Range("K67").Value = 7167
$ws.Range("M67").Value = -6231
$ws.Range("H122").Value = 1482.5652
$ws.Range("I122").Value = 1319.375
$ws.Range("J122").Value = 1569.6
$ws.Range("K122").Value = 11874.375
$ws.Range("L122").Value = 14126.4
$ws.Range("M122").Value = -9424.375
$ws.Range("N122").Value = -19026.4
$ws.Range("H129").Value = 2529.6667
$ws.Range("I129").Value = 2005
$ws.Range("J129").Value = 2792
$ws.Range("K129").Value = 6015
$ws.Range("L129").Value = 8376
$ws.Range("M129").Value = -1015
$ws.Range("N129").Value = -18376

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 483.61905
$ws.Range("I97").Value = 426.27274
$ws.Range("K97").Value = 426.27274
$ws.Range("M97").Value = 69.72726

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 4400
$ws.Range("I55").Value = 4400
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 4400
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("M55").Value = -4227
$ws.Range("H93").Value = 1434999.9
$ws.Range("I93").Value = 8250
$ws.Range("J93").Value = 3337333
$ws.Range("K93").Value = 8250
$ws.Range("L93").Value = 3337333
$ws.Range("M93").Value = -7002
$ws.Range("N93").Value = -3339829

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 502999.5
$ws.Range("J4").Value = 6000
$ws.Range("L4").Value = 6000
$ws.Range("N4").Value = -6226
$ws.Range("H96").Value = 2924
$ws.Range("I96").Value = 2815.5
$ws.Range("J96").Value = 2996.3333
$ws.Range("K96").Value = 2815.5
$ws.Range("L96").Value = 2996.3333
$ws.Range("M96").Value = -1442.5
$ws.Range("N96").Value = -5742.3333
$ws.Range("H100").Value = 754.4545000000001
$ws.Range("I100").Value = 765.1579
$ws.Range("J100").Value = 686.6667
$ws.Range("K100").Value = 1530.3158
$ws.Range("L100").Value = 1373.3334
$ws.Range("M100").Value = -989.3158000000001
$ws.Range("N100").Value = -2455.3334
$ws.Range("H136").Value = 2384.7805
$ws.Range("I136").Value = 1271.8064
$ws.Range("J136").Value = 5835
$ws.Range("K136").Value = 3815.4192
$ws.Range("L136").Value = 17505
$ws.Range("M136").Value = -1265.4192
$ws.Range("N136").Value = -22605
